$d = $word.ActiveDocument

# Locate the paragraph that starts the "shared folder" block so we can
# insert the new content right before it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("The shared folder is stored in")) {
        $target = $p
        break
    }
}

# The existing (hidden) _GoBack bookmark currently sits at the end of the
# "Created a symbolic link ... to /home/vagrant" paragraph. It is being
# retyped as part of this edit, so drop it here - it gets re-created at
# its new location (inside the freshly typed paragraph) by the inserted
# XML below.
$bms = $d.Bookmarks
if ($bms.Exists("_GoBack")) {
    $bms.Item("_GoBack").Delete()
}

# Replace "The shared folder is stored in " paragraph with: a new blank
# paragraph, a new paragraph describing the python script (with the
# _GoBack bookmark sitting where the cursor was left after typing), and
# then the original "The shared folder is stored in " paragraph text
# (so its content is preserved, just re-emitted).
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r><w:t xml:space="preserve">The </w:t></w:r>
            <w:r><w:t xml:space="preserve">python </w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t>script writes to the folder /home/vagrant/Mail folder</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:t xml:space="preserve">The shared folder is stored in </w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($xml)
